$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated symbol list values from coinranking.com scrape
$updates = @(
    @{ Cell = 'D2'; Value = '245.23' }
    @{ Cell = 'E2'; Value = '-0.87%' }
    @{ Cell = 'D3'; Value = '27.32' }
    @{ Cell = 'E3'; Value = '4.05%' }
    @{ Cell = 'D4'; Value = '5.119' }
    @{ Cell = 'E4'; Value = '0.55%' }
    @{ Cell = 'D5'; Value = '0.05684' }
    @{ Cell = 'E5'; Value = '1.49%' }
    @{ Cell = 'D6'; Value = '6.524' }
    @{ Cell = 'E6'; Value = '0.75%' }
    @{ Cell = 'D7'; Value = '0.8201' }
    @{ Cell = 'E7'; Value = '0.91%' }
    @{ Cell = 'D8'; Value = '0.8610' }
    @{ Cell = 'E8'; Value = '1.93%' }
    @{ Cell = 'D9'; Value = '0.1334' }
    @{ Cell = 'E9'; Value = '-0.45%' }
    @{ Cell = 'D10'; Value = '0.06939' }
    @{ Cell = 'E10'; Value = '-0.75%' }
    @{ Cell = 'D11'; Value = '0.02860' }
    @{ Cell = 'E11'; Value = '1.88%' }
    @{ Cell = 'E12'; Value = '0.02%' }
    @{ Cell = 'E13'; Value = '0.51%' }
    @{ Cell = 'E14'; Value = '-12.39%' }
    @{ Cell = 'B15'; Value = 'TigerCash' }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' }
    @{ Cell = 'D15'; Value = '0.006213' }
    @{ Cell = 'E15'; Value = '0.79%' }
    @{ Cell = 'B16'; Value = 'LEO' }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' }
    @{ Cell = 'D16'; Value = '3.511' }
    @{ Cell = 'E16'; Value = '-2.62%' }
    @{ Cell = 'B17'; Value = 'GateToken' }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' }
    @{ Cell = 'D17'; Value = '3.010' }
    @{ Cell = 'E17'; Value = '-0.30%' }
    @{ Cell = 'B18'; Value = 'BTSEToken' }
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse' }
    @{ Cell = 'D18'; Value = '2.313' }
    @{ Cell = 'E18'; Value = '12.56%' }
    @{ Cell = 'B19'; Value = 'One' }
    @{ Cell = 'C19'; Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one' }
    @{ Cell = 'D19'; Value = '0.01017' }
    @{ Cell = 'E19'; Value = '1,605.94%' }
    @{ Cell = 'D20'; Value = '0.3145' }
    @{ Cell = 'E20'; Value = '1.15%' }
    @{ Cell = 'D21'; Value = '0.03206' }
    @{ Cell = 'E21'; Value = '0.40%' }
    @{ Cell = 'E22'; Value = '-0.05%' }
    @{ Cell = 'D23'; Value = '3.550' }
    @{ Cell = 'E23'; Value = '-5.69%' }
    @{ Cell = 'E24'; Value = '1.75%' }
    @{ Cell = 'D25'; Value = '0.001218' }
    @{ Cell = 'E25'; Value = '-2.06%' }
    @{ Cell = 'D26'; Value = '0.004467' }
    @{ Cell = 'E26'; Value = '-2.40%' }
    @{ Cell = 'E27'; Value = '22.88%' }
    @{ Cell = 'D40'; Value = '0.03718' }
    @{ Cell = 'E40'; Value = '1.51%' }
    @{ Cell = 'D41'; Value = '0.005958' }
    @{ Cell = 'E41'; Value = '-3.59%' }
    @{ Cell = 'E42'; Value = '0.23%' }
    @{ Cell = 'D43'; Value = '0.002299' }
    @{ Cell = 'E43'; Value = '-8.03%' }
    @{ Cell = 'D44'; Value = '0.009712' }
    @{ Cell = 'E44'; Value = '17.34%' }
    @{ Cell = 'D45'; Value = '0.00005114' }
    @{ Cell = 'E45'; Value = '-5.07%' }
    @{ Cell = 'E46'; Value = '-0.03%' }
    @{ Cell = 'E48'; Value = '3.18%' }
    @{ Cell = 'D49'; Value = '0.00002099' }
    @{ Cell = 'E49'; Value = '-0.03%' }
    @{ Cell = 'D50'; Value = '0.0001999' }
    @{ Cell = 'E50'; Value = '-0.03%' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}

Write-Host "Applied $($updates.Count) cell updates"